# Reposition/resize the two "sandboxed" shapes on the VEE diagram slide so
# the apps sit in-between the VEE diagram (highlighting OTA).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Target EMU values (sub-point precision is chosen to land on the exact
# EMU after the COM layer's internal point<->EMU rounding):
#   shape1: off.x 3968151 (unchanged), off.y 1839757, ext.cy 1835208
#   shape2: off.x 4037153, off.y 1839758, ext.cy 363609

# Shape "Rounded Rectangle 105" (the accent-colored sandbox body)
$shp1 = $s.Shapes.Item(15)
$shp1.Top    = 144.8627929805118
$shp1.Height = 144.50460815413385

# Shape "Round Same Side Corner Rectangle 106" (the rotated header band)
$shp2 = $s.Shapes.Item(16)
$shp2.Left   = 317.88610839212595
$shp2.Top    = 144.8628692706693
$shp2.Height = 28.630668646259842
